$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 5: Number=3, Tool=6, Borrower=Trunks, Date, Expiration
Set-TextValue $ws.Cells.Item(5,1) "3"
Set-TextValue $ws.Cells.Item(5,2) "6"
$ws.Cells.Item(5,3).Value = "Trunks"
$ws.Cells.Item(5,4).Value = "Wed May 23 16:12:00 GMT-03:00 2018"
$ws.Cells.Item(5,5).Value = "Thu May 24 16:12:00 GMT-03:00 2018"

# Row 6: Number=4, Tool=3, Borrower=Vegeta, Date, Expiration
Set-TextValue $ws.Cells.Item(6,1) "4"
Set-TextValue $ws.Cells.Item(6,2) "3"
$ws.Cells.Item(6,3).Value = "Vegeta"
$ws.Cells.Item(6,4).Value = "Wed May 23 16:12:14 GMT-03:00 2018"
$ws.Cells.Item(6,5).Value = "Thu May 24 16:12:14 GMT-03:00 2018"

# Row 7: Number=5, Tool=2, Borrower=Gohan, Date, Expiration
Set-TextValue $ws.Cells.Item(7,1) "5"
Set-TextValue $ws.Cells.Item(7,2) "2"
$ws.Cells.Item(7,3).Value = "Gohan"
$ws.Cells.Item(7,4).Value = "Wed May 23 16:12:36 GMT-03:00 2018"
$ws.Cells.Item(7,5).Value = "Thu May 24 16:12:36 GMT-03:00 2018"

# Row 8: Number=6, Tool=4, Borrower=Gohan, Date, Expiration
Set-TextValue $ws.Cells.Item(8,1) "6"
Set-TextValue $ws.Cells.Item(8,2) "4"
$ws.Cells.Item(8,3).Value = "Gohan"
$ws.Cells.Item(8,4).Value = "Wed May 23 16:12:42 GMT-03:00 2018"
$ws.Cells.Item(8,5).Value = "Thu May 24 16:12:42 GMT-03:00 2018"
